$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.006697333333333
$ws.Range("H2").Value = 3.020092
$ws.Range("I2").Value = 0.0001985651645046208
$ws.Range("J2").Value = 0.0001985651645046208
$ws.Range("M2").Value = 0.1809866666666667
$ws.Range("N2").Value = 0.54296
$ws.Range("O2").Value = 0.03987407676082905
$ws.Range("P2").Value = 0.03987407676082905
$ws.Range("Q2").Value = 0.1821987947022222
$ws.Range("R2").Value = 1.63978915232
$ws.Range("S2").Value = 0.000007917602611483899
$ws.Range("T2").Value = 0.000007917602611483899
$ws.Range("G3").Value = 1.006697333333333
$ws.Range("H3").Value = 3.020092
$ws.Range("I3").Value = 0.0001985651645046208
$ws.Range("J3").Value = 0.0001985651645046208
$ws.Range("O3").Value = 0.1057193993302571
$ws.Range("P3").Value = 0.1057193993302571
$ws.Range("Q3").Value = 0.4830694200182222
$ws.Range("R3").Value = 4.347624780164
$ws.Range("S3").Value = 0.0000209921899193422
$ws.Range("T3").Value = 0.0000209921899193422
$ws.Range("G4").Value = 1.006697333333333
$ws.Range("H4").Value = 3.020092
$ws.Range("I4").Value = 0.0001985651645046208
$ws.Range("J4").Value = 0.0001985651645046208
$ws.Range("M4").Value = 3.878113333333333
$ws.Range("N4").Value = 11.63434
$ws.Range("O4").Value = 0.8544065239089139
$ws.Range("P4").Value = 0.8544065239089139
$ws.Range("Q4").Value = 3.904086351031111
$ws.Range("R4").Value = 35.13677715928
$ws.Range("S4").Value = 0.0001696553719737947
$ws.Range("T4").Value = 0.0001696553719737947
$ws.Range("I5").Value = 0.9806494927176636
$ws.Range("J5").Value = 0.9806494927176637
$ws.Range("M5").Value = 0.1809866666666667
$ws.Range("N5").Value = 0.54296
$ws.Range("O5").Value = 0.03987407676082905
$ws.Range("P5").Value = 0.03987407676082905
$ws.Range("Q5").Value = 899.821255376071
$ws.Range("R5").Value = 8098.391298384639
$ws.Range("S5").Value = 0.03910249314809219
$ws.Range("T5").Value = 0.03910249314809219
$ws.Range("I6").Value = 0.9806494927176636
$ws.Range("J6").Value = 0.9806494927176637
$ws.Range("O6").Value = 0.1057193993302571
$ws.Range("P6").Value = 0.1057193993302571
$ws.Range("S6").Value = 0.1036736753236327
$ws.Range("T6").Value = 0.1036736753236327
$ws.Range("I7").Value = 0.9806494927176636
$ws.Range("J7").Value = 0.9806494927176637
$ws.Range("M7").Value = 3.878113333333333
$ws.Range("N7").Value = 11.63434
$ws.Range("O7").Value = 0.8544065239089139
$ws.Range("P7").Value = 0.8544065239089139
$ws.Range("Q7").Value = 19281.02700801539
$ws.Range("R7").Value = 173529.2430721385
$ws.Range("S7").Value = 0.8378733242459387
$ws.Range("T7").Value = 0.8378733242459389
$ws.Range("G8").Value = 97.097641
$ws.Range("H8").Value = 291.292923
$ws.Range("I8").Value = 0.01915194211783179
$ws.Range("J8").Value = 0.01915194211783179
$ws.Range("M8").Value = 0.1809866666666667
$ws.Range("N8").Value = 0.54296
$ws.Range("O8").Value = 0.03987407676082905
$ws.Range("P8").Value = 0.03987407676082905
$ws.Range("Q8").Value = 17.57337838578666
$ws.Range("R8").Value = 158.16040547208
$ws.Range("S8").Value = 0.0007636660101253797
$ws.Range("T8").Value = 0.0007636660101253797
$ws.Range("G9").Value = 97.097641
$ws.Range("H9").Value = 291.292923
$ws.Range("I9").Value = 0.01915194211783179
$ws.Range("J9").Value = 0.01915194211783179
$ws.Range("O9").Value = 0.1057193993302571
$ws.Range("P9").Value = 0.1057193993302571
$ws.Range("Q9").Value = 46.59285325381567
$ws.Range("R9").Value = 419.335679284341
$ws.Range("S9").Value = 0.002024731816705029
$ws.Range("T9").Value = 0.002024731816705029
$ws.Range("G10").Value = 97.097641
$ws.Range("H10").Value = 291.292923
$ws.Range("I10").Value = 0.01915194211783179
$ws.Range("J10").Value = 0.01915194211783179
$ws.Range("M10").Value = 3.878113333333333
$ws.Range("N10").Value = 11.63434
$ws.Range("O10").Value = 0.8544065239089139
$ws.Range("P10").Value = 0.8544065239089139
$ws.Range("Q10").Value = 376.5556561973133
$ws.Range("R10").Value = 3389.00090577582
$ws.Range("S10").Value = 0.01636354429100138
$ws.Range("T10").Value = 0.01636354429100138
